# Daily attendance processing - reorder "Recorded By" (column G) entries so
# that "System" / "system" is listed first among the comma-separated names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value()

    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
    elseif ($val -eq "system, System, backup@backdoor.com") {
        $cell.Value = "System, system, backup@backdoor.com"
    }
}
